$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.270298957824707
$ws.Range("B1").Value = 2.150876998901367
$ws.Range("C1").Value = 4.676649570465088
$ws.Range("D1").Value = 3.123369693756104
$ws.Range("E1").Value = 1.375371336936951
